$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 14603.857
$ws.Range("J9").Value = 800
$ws.Range("L9").Value = 800
$ws.Range("N9").Value = -1138

$ws.Range("H18").Value = 723.5833
$ws.Range("I18").Value = 698.4545000000001
$ws.Range("K18").Value = 698.4545000000001
$ws.Range("M18").Value = -414.4545000000001

$ws.Range("H28").Value = 1019.0909
$ws.Range("I28").Value = 919.125
$ws.Range("K28").Value = 919.125
$ws.Range("M28").Value = -434.125

$ws.Range("H33").Value = 232
$ws.Range("I33").Value = 213.4
$ws.Range("J33").Value = 325
$ws.Range("K33").Value = 213.4
$ws.Range("L33").Value = 325
$ws.Range("M33").Value = 15.59999999999999
$ws.Range("N33").Value = -783

$ws.Range("H55").Value = 56
$ws.Range("I55").Value = 56
$ws.Range("K55").Value = 56
$ws.Range("M55").Value = 158

$ws.Range("H61").Value = 820
$ws.Range("I61").Value = 820
$ws.Range("K61").Value = 2460
$ws.Range("M61").Value = -2288

$ws.Range("H70").Value = 3743.9546
$ws.Range("I70").Value = 1272.3334
$ws.Range("J70").Value = 6709.9
$ws.Range("K70").Value = 3817.0002
$ws.Range("L70").Value = 20129.7
$ws.Range("M70").Value = -3547.0002
$ws.Range("N70").Value = -20669.7

$ws.Range("H73").Value = 3743.9546
$ws.Range("I73").Value = 1272.3334
$ws.Range("J73").Value = 6709.9
$ws.Range("K73").Value = 3817.0002
$ws.Range("L73").Value = 20129.7
$ws.Range("M73").Value = -2881.0002
$ws.Range("N73").Value = -22001.7

$ws.Range("H98").Value = 1204.8462
$ws.Range("I98").Value = 884.4167
$ws.Range("K98").Value = 884.4167
$ws.Range("M98").Value = 613.5833

$ws.Range("H113").Value = 3834.6667
$ws.Range("I113").Value = 2002
$ws.Range("K113").Value = 2002
$ws.Range("M113").Value = 1252

$ws.Range("H122").Value = 1204.8462
$ws.Range("I122").Value = 884.4167
$ws.Range("K122").Value = 2653.2501
$ws.Range("M122").Value = -203.2501000000002

$ws.Range("H132").Value = 3236.125
$ws.Range("I132").Value = 3257.9092
$ws.Range("K132").Value = 9773.7276
$ws.Range("M132").Value = -7243.7276

$ws.Range("H137").Value = 2412.5483
$ws.Range("I137").Value = 2089
$ws.Range("J137").Value = 3342.75
$ws.Range("K137").Value = 6267
$ws.Range("L137").Value = 10028.25
$ws.Range("M137").Value = -3717
$ws.Range("N137").Value = -15128.25

$ws.Range("H141").Value = 1738.6666
$ws.Range("I141").Value = 1738.6666
$ws.Range("K141").Value = 5215.9998
$ws.Range("M141").Value = -35.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5670.375
$ws.Range("J45").Value = 2276.75
$ws.Range("L45").Value = 2276.75
$ws.Range("N45").Value = -3030.75

$ws.Range("H63").Value = 1999.5
$ws.Range("J63").Value = 1999.3334
$ws.Range("L63").Value = 1999.3334
$ws.Range("N63").Value = -3371.3334

$ws.Range("H66").Value = 1999.5
$ws.Range("J66").Value = 1999.3334
$ws.Range("L66").Value = 9996.666999999999
$ws.Range("N66").Value = -16860.667

$ws.Range("H122").Value = 4045.487
$ws.Range("I122").Value = 3461.75
$ws.Range("J122").Value = 6714
$ws.Range("K122").Value = 10385.25
$ws.Range("L122").Value = 20142
$ws.Range("M122").Value = -7935.25
$ws.Range("N122").Value = -25042

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2282.75
$ws.Range("I22").Value = 2983
$ws.Range("J22").Value = 1582.5
$ws.Range("K22").Value = 2983
$ws.Range("L22").Value = 1582.5
$ws.Range("M22").Value = -2810
$ws.Range("N22").Value = -1928.5

$ws.Range("H64").Value = 720
$ws.Range("I64").Value = 533
$ws.Range("J64").Value = 782.3333
$ws.Range("K64").Value = 533
$ws.Range("L64").Value = 782.3333
$ws.Range("M64").Value = -308
$ws.Range("N64").Value = -1232.3333

$ws.Range("H67").Value = 720
$ws.Range("I67").Value = 533
$ws.Range("J67").Value = 782.3333
$ws.Range("K67").Value = 533
$ws.Range("L67").Value = 782.3333
$ws.Range("M67").Value = 247
$ws.Range("N67").Value = -2342.3333

$ws.Range("H94").Value = 3115.889
$ws.Range("I94").Value = 3335
$ws.Range("K94").Value = 3335
$ws.Range("M94").Value = -2884

$ws.Range("H99").Value = 2009.6666
$ws.Range("J99").Value = 2009.6666
$ws.Range("L99").Value = 2009.6666
$ws.Range("N99").Value = -5005.6666

$ws.Range("H134").Value = 23813610
$ws.Range("I134").Value = 27781774
$ws.Range("J134").Value = 4632
$ws.Range("K134").Value = 83345322
$ws.Range("L134").Value = 13896
$ws.Range("M134").Value = -83342787
$ws.Range("N134").Value = -18966

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 447.5
$ws.Range("I2").Value = 447.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 447.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -334.5
$ws.Range("N2").ClearContents()

$ws.Range("H7").Value = 1302.6923
$ws.Range("J7").Value = 73
$ws.Range("L7").Value = 73
$ws.Range("N7").Value = -299

$ws.Range("H9").Value = 32999.5
$ws.Range("J9").Value = 32999.5
$ws.Range("L9").Value = 32999.5
$ws.Range("N9").Value = -33335.5

$ws.Range("H16").Value = 1098273.5
$ws.Range("I16").Value = 1371888.8
$ws.Range("J16").Value = 3812.5
$ws.Range("K16").Value = 1371888.8
$ws.Range("L16").Value = 3812.5
$ws.Range("M16").Value = -1371601.8
$ws.Range("N16").Value = -4386.5

$ws.Range("H31").Value = 4343.25
$ws.Range("I31").Value = 4999.3335
$ws.Range("J31").Value = 2375
$ws.Range("K31").Value = 4999.3335
$ws.Range("L31").Value = 2375
$ws.Range("M31").Value = -4704.3335
$ws.Range("N31").Value = -2965

$ws.Range("H34").Value = 4343.25
$ws.Range("I34").Value = 4999.3335
$ws.Range("J34").Value = 2375
$ws.Range("K34").Value = 4999.3335
$ws.Range("L34").Value = 2375
$ws.Range("M34").Value = -4797.3335
$ws.Range("N34").Value = -2779

$ws.Range("H58").Value = 62513164
$ws.Range("I58").Value = 62513164
$ws.Range("K58").Value = 62513164
$ws.Range("M58").Value = -62512961

$ws.Range("H98").Value = 88997.5
$ws.Range("J98").Value = 88997.5
$ws.Range("L98").Value = 88997.5
$ws.Range("N98").Value = -93489.5

$ws.Range("H113").Value = 1098273.5
$ws.Range("I113").Value = 1371888.8
$ws.Range("J113").Value = 3812.5
$ws.Range("K113").Value = 1371888.8
$ws.Range("L113").Value = 3812.5
$ws.Range("M113").Value = -1369718.8
$ws.Range("N113").Value = -8152.5

$ws.Range("H122").Value = 3994.1667
$ws.Range("I122").Value = 2630
$ws.Range("J122").Value = 5358.3335
$ws.Range("K122").Value = 7890
$ws.Range("L122").Value = 16075.0005
$ws.Range("M122").Value = -5440
$ws.Range("N122").Value = -20975.0005

$ws.Range("H134").Value = 5841467
$ws.Range("I134").Value = 6440034.5
$ws.Range("K134").Value = 19320103.5
$ws.Range("M134").Value = -19317568.5

$ws.Range("H136").Value = 62513164
$ws.Range("I136").Value = 62513164
$ws.Range("K136").Value = 187539492
$ws.Range("M136").Value = -187536942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 534.7273
$ws.Range("I2").Value = 36.4
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 218.4
$ws.Range("L2").Value = 5700
$ws.Range("M2").Value = -105.4
$ws.Range("N2").Value = -5926

$ws.Range("H7").Value = 1432359.8
$ws.Range("I7").Value = 2502749.8
$ws.Range("J7").Value = 5173
$ws.Range("K7").Value = 7508249.399999999
$ws.Range("L7").Value = 15519
$ws.Range("M7").Value = -7508137.399999999
$ws.Range("N7").Value = -15743

$ws.Range("H10").Value = 153.33333
$ws.Range("J10").Value = 40
$ws.Range("L10").Value = 120
$ws.Range("N10").Value = -398

$ws.Range("H11").Value = 198964.22
$ws.Range("J11").Value = 66833
$ws.Range("L11").Value = 200499
$ws.Range("N11").Value = -200779

$ws.Range("H12").Value = 175
$ws.Range("J12").Value = 268.1111
$ws.Range("L12").Value = 804.3333
$ws.Range("N12").Value = -1150.3333

$ws.Range("H13").Value = 347.14285
$ws.Range("I13").Value = 330
$ws.Range("J13").Value = 450
$ws.Range("K13").Value = 990
$ws.Range("L13").Value = 1350
$ws.Range("M13").Value = -822
$ws.Range("N13").Value = -1686

$ws.Range("H22").Value = 4001
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4001
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 12003
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -12341

$ws.Range("H23").Value = 403.8
$ws.Range("J23").Value = 477.25
$ws.Range("L23").Value = 1431.75
$ws.Range("N23").Value = -1901.75

$ws.Range("H25").Value = 1275
$ws.Range("I25").Value = 800
$ws.Range("K25").Value = 2400
$ws.Range("M25").Value = -2231

$ws.Range("H27").Value = 4001
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4001
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 12003
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -12207

$ws.Range("H30").Value = 1275
$ws.Range("I30").Value = 800
$ws.Range("K30").Value = 2400
$ws.Range("M30").Value = -2298

$ws.Range("H34").Value = 1232.6666
$ws.Range("J34").Value = 1799
$ws.Range("L34").Value = 5397
$ws.Range("N34").Value = -5565

$ws.Range("H37").Value = 116026.53
$ws.Range("J37").Value = 116026.53
$ws.Range("L37").Value = 348079.59
$ws.Range("N37").Value = -348303.59

$ws.Range("H39").Value = 958.9167
$ws.Range("I39").Value = 774.2727
$ws.Range("J39").Value = 2990
$ws.Range("K39").Value = 2322.8181
$ws.Range("L39").Value = 8970
$ws.Range("M39").Value = -2028.8181
$ws.Range("N39").Value = -9558

$ws.Range("H44").Value = 1666.3334
$ws.Range("I44").Value = 1249.75
$ws.Range("J44").Value = 2499.5
$ws.Range("K44").Value = 3749.25
$ws.Range("L44").Value = 7498.5
$ws.Range("M44").Value = -3351.25
$ws.Range("N44").Value = -8294.5

$ws.Range("H46").Value = 302
$ws.Range("J46").Value = 302
$ws.Range("L46").Value = 906
$ws.Range("N46").Value = -1088

$ws.Range("H47").Value = 307.6
$ws.Range("I47").Value = 283.5
$ws.Range("J47").Value = 404
$ws.Range("K47").Value = 850.5
$ws.Range("L47").Value = 1212
$ws.Range("M47").Value = -419.5
$ws.Range("N47").Value = -2074

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H64").Value = 2599.6667
$ws.Range("I64").Value = 1399.5
$ws.Range("K64").Value = 4198.5
$ws.Range("M64").Value = -3928.5

$ws.Range("H67").Value = 2599.6667
$ws.Range("I67").Value = 1399.5
$ws.Range("K67").Value = 4198.5
$ws.Range("M67").Value = -3262.5

$ws.Range("H68").Value = 1501.7
$ws.Range("J68").Value = 2403.8
$ws.Range("L68").Value = 7211.400000000001
$ws.Range("N68").Value = -8833.400000000001

$ws.Range("H69").Value = 648.2
$ws.Range("I69").Value = 648.2
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 1944.6
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -1133.6
$ws.Range("N69").ClearContents()

$ws.Range("H71").Value = 1501.7
$ws.Range("J71").Value = 2403.8
$ws.Range("L71").Value = 21634.2
$ws.Range("N71").Value = -29746.2

$ws.Range("H72").Value = 648.2
$ws.Range("I72").Value = 648.2
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 5833.8
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -1777.8
$ws.Range("N72").ClearContents()

$ws.Range("H76").Value = 19998.166
$ws.Range("I76").Value = 19994.5
$ws.Range("K76").Value = 59983.5
$ws.Range("M76").Value = -59600.5

$ws.Range("H79").Value = 19998.166
$ws.Range("I79").Value = 19994.5
$ws.Range("K79").Value = 59983.5
$ws.Range("M79").Value = -58657.5

$ws.Range("H80").Value = 2050
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2050
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6150
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -8022

$ws.Range("H83").Value = 2050
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2050
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18450
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -27810

$ws.Range("H107").Value = 664.6875
$ws.Range("J107").Value = 999.625
$ws.Range("L107").Value = 2998.875
$ws.Range("N107").Value = -6838.875

$ws.Range("H121").Value = 102198.6
$ws.Range("I121").Value = 167393.17
$ws.Range("J121").Value = 4406.75
$ws.Range("K121").Value = 502179.51
$ws.Range("L121").Value = 13220.25
$ws.Range("M121").Value = -500869.51
$ws.Range("N121").Value = -15840.25

$ws.Range("H134").Value = 1229.75
$ws.Range("I134").Value = 1229.75
$ws.Range("K134").Value = 3689.25
$ws.Range("M134").Value = 1380.75

$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -2400

$ws.Range("H140").Value = 497.53845
$ws.Range("I140").Value = 497.53845
$ws.Range("K140").Value = 1492.61535
$ws.Range("M140").Value = 3687.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1549.5897
$ws.Range("I97").Value = 1445.88
$ws.Range("J97").Value = 1734.7858
$ws.Range("K97").Value = 1445.88
$ws.Range("L97").Value = 1734.7858
$ws.Range("M97").Value = -949.8800000000001
$ws.Range("N97").Value = -2726.7858

$ws.Range("H132").Value = 5212278
$ws.Range("I132").Value = 6253671
$ws.Range("K132").Value = 18761013
$ws.Range("M132").Value = -18758483

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3841.9285
$ws.Range("I7").Value = 3599.6
$ws.Range("J7").Value = 4447.75
$ws.Range("K7").Value = 3599.6
$ws.Range("L7").Value = 4447.75
$ws.Range("M7").Value = -3487.6
$ws.Range("N7").Value = -4671.75

$ws.Range("H16").Value = 1636.2
$ws.Range("I16").Value = 1191.2142
$ws.Range("J16").Value = 2674.5
$ws.Range("K16").Value = 1191.2142
$ws.Range("L16").Value = 2674.5
$ws.Range("M16").Value = -1021.2142
$ws.Range("N16").Value = -3014.5

$ws.Range("H46").Value = 1700.5
$ws.Range("I46").Value = 1767.2222
$ws.Range("K46").Value = 1767.2222
$ws.Range("M46").Value = -1579.2222

$ws.Range("H68").Value = 13161894
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 13161894
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H93").Value = 2837.5557
$ws.Range("I93").Value = 999.6
$ws.Range("K93").Value = 999.6
$ws.Range("M93").Value = 248.4

$ws.Range("H100").Value = 9981667
$ws.Range("I100").Value = 24949950
$ws.Range("K100").Value = 24949950
$ws.Range("M100").Value = -24949409

$ws.Range("H104").Value = 12999
$ws.Range("J104").Value = 12999
$ws.Range("L104").Value = 12999
$ws.Range("N104").Value = -19987

$ws.Range("H126").Value = 3841.9285
$ws.Range("I126").Value = 3599.6
$ws.Range("J126").Value = 4447.75
$ws.Range("K126").Value = 10798.8
$ws.Range("L126").Value = 13343.25
$ws.Range("M126").Value = -8328.799999999999
$ws.Range("N126").Value = -18283.25

$ws.Range("H132").Value = 5440812
$ws.Range("I132").Value = 5958408.5
$ws.Range("J132").Value = 6047.25
$ws.Range("K132").Value = 17875225.5
$ws.Range("L132").Value = 18141.75
$ws.Range("M132").Value = -17872695.5
$ws.Range("N132").Value = -23201.75

$ws.Range("H136").Value = 1783.1482
$ws.Range("I136").Value = 1515.6471
$ws.Range("J136").Value = 2237.9
$ws.Range("K136").Value = 4546.9413
$ws.Range("L136").Value = 6713.700000000001
$ws.Range("M136").Value = -1996.9413
$ws.Range("N136").Value = -11813.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H62").Value = 6750
$ws.Range("I62").Value = 6750
$ws.Range("K62").Value = 6750
$ws.Range("M62").Value = -6126

$ws.Range("H65").Value = 6750
$ws.Range("I65").Value = 6750
$ws.Range("K65").Value = 33750
$ws.Range("M65").Value = -30630

$ws.Range("H113").Value = 761.5
$ws.Range("I113").Value = 732.64
$ws.Range("K113").Value = 2197.92
$ws.Range("M113").Value = -27.92000000000007

$ws.Range("H127").Value = 70000
$ws.Range("J127").Value = 70000
$ws.Range("L127").Value = 70000
$ws.Range("N127").Value = -79920

$ws.Range("H132").Value = 9348658
$ws.Range("I132").Value = 11495954
$ws.Range("K132").Value = 34487862
$ws.Range("M132").Value = -34485332

$ws.Range("H136").Value = 27779904
$ws.Range("I136").Value = 29413722
$ws.Range("K136").Value = 88241166
$ws.Range("M136").Value = -88238616
